$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to remain text (values look numeric)
foreach ($addr in @("D5","D6","D8","D9","D10","D11","D15","D16","D18","D20","D22","D23","D24","D26","D29","D30","D31","D37","D39","D40","D43","D46")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "27.514.11"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.621.83"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "211.59"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "23.11"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "0.264"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").Value = "0.0611"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "0.0883"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").Value = "1.850.82"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "1.615.89"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "0.550"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "65.41"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "27.489.99"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "229.83"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "10.45"
$ws.Range("E22").Value = "  +4.32%  "
$ws.Range("D23").Value = "4.35"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +8.95%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "6.88"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").Value = "1.465.18"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "0.952"
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").Value = "0.873"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "0.552"
$ws.Range("E40").Value = "  -2.75%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").Value = "67.30"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.760.90"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  +0.20%  "

Write-Host "Applied cryptos update"
